$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AS (column index 45).
# This shifts the old "nom" column (AS) to AT, and old "url_produit" (AT) to AU.
$ws.Columns.Item(45).Insert()

# The new AS column represents a fresh price snapshot taken at 2026-01-29 16:28:56.
$ws.Range("AS1").Value = "2026-01-29 16:28:56"

# Find the last used row so we copy the latest price (column AR) into the new
# snapshot column (AS) for every product row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 44).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $prevValue = $ws.Cells.Item($row, 44).Value()
    $ws.Cells.Item($row, 45).Value = $prevValue
}
